$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(413).Insert()

$ws.Range("A413").Value = 4
$ws.Range("B413").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C413").Value = "Los Lagos"
$ws.Range("D413").Value = 45267
$ws.Range("E413").Value = 10
$ws.Range("F413").Value = 100112039
$ws.Range("G413").Value = "Ciboulette"
$ws.Range("H413").Value = "Sin especificar"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 120
$ws.Range("K413").Value = 2500
$ws.Range("L413").Value = 2500
$ws.Range("M413").Value = 2500
$ws.Range("N413").Value = "$/docena de atados"
$ws.Range("O413").Value = "Región Metropolitana"
$ws.Range("P413").Value = 833
$ws.Range("Q413").Value = 3
$ws.Range("R413").Value = "Hortaliza"
